# Implements "stratification into the tool":
#   - The single "fragility_curve" bucket used across loads / generators /
#     transformers / lines (previously "dummy" / "tower_water") is replaced
#     by two explicit strata: "towers_1" and "towers_2".
#   - loads / generators / transformers gain a new "normalTTR" column set to 1.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# loads sheet: drop the old fragility_curve value in column P and add a
# new normalTTR column (Q) with value 1 for every data row.
# ---------------------------------------------------------------------
$loads = $wb.Worksheets.Item("loads")
$loads.Range("Q1").Value = "normalTTR"
for ($r = 2; $r -le 9; $r++) {
    $loads.Cells.Item($r, 16).Value = ""   # column P
    $loads.Cells.Item($r, 17).Value = 1    # column Q
}

# ---------------------------------------------------------------------
# generators sheet: drop the old fragility_curve value in column W and add
# a new normalTTR column (Z) with value 1 for every data row.
# ---------------------------------------------------------------------
$generators = $wb.Worksheets.Item("generators")
$generators.Range("Z1").Value = "normalTTR"
for ($r = 2; $r -le 7; $r++) {
    $generators.Cells.Item($r, 23).Value = ""   # column W
    $generators.Cells.Item($r, 26).Value = 1    # column Z
}

# ---------------------------------------------------------------------
# transformers sheet: drop the old fragility_curve value in column E and
# add a new normalTTR column (L) with value 1.
# ---------------------------------------------------------------------
$transformers = $wb.Worksheets.Item("transformers")
$transformers.Range("L1").Value = "normalTTR"
$transformers.Range("E2").Value = ""
$transformers.Range("L2").Value = 1

# ---------------------------------------------------------------------
# lines sheet: column K held the single "tower_water" fragility_curve for
# every line; it is now split between the two new strata, and a few lines
# no longer reference any fragility curve at all.
# ---------------------------------------------------------------------
$lines = $wb.Worksheets.Item("lines")
$linesK = @{
    2  = ""
    3  = "towers_1"
    4  = "towers_2"
    5  = "towers_1"
    6  = "towers_1"
    7  = "towers_2"
    8  = "towers_2"
    9  = "towers_2"
    10 = "towers_1"
    11 = ""
    12 = "towers_1"
    13 = "towers_1"
    14 = "towers_2"
    15 = "towers_2"
    16 = ""
    17 = ""
}
foreach ($r in $linesK.Keys) {
    $lines.Cells.Item($r, 11).Value = $linesK[$r]   # column K
}
